# Update countries & provincias Spain
# Applies the data refresh captured in the diff: updated case counts for
# several countries (which also causes a few adjacent countries to swap
# rank/position in the table), and updates the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" banner text (row 1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 03:22"

# Helper: write a full data row (Pais, Casos totales, Nuevos casos,
# Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
function Set-Row($Row, $Values) {
    $arr = New-Object 'object[,]' 1,8
    for ($i = 0; $i -lt 8; $i++) {
        $arr[0,$i] = $Values[$i]
    }
    $addr = "A" + $Row + ":H" + $Row
    $ws.Range($addr).Value = $arr
}

# Estados Unidos
Set-Row 4 @("Estados Unidos", 763832, 0, 71003, 652276, 13566, 0, 40553)

# China
Set-Row 11 @("China", 82747, 12, 77084, 1031, 81, 0, 4632)

# Brasil
Set-Row 14 @("Brasil", 38654, 0, 14026, 22166, 6634, 0, 2462)

# Corea del Sur
Set-Row 27 @("Corea del Sur", 10674, 13, 8114, 2324, 55, 2, 236)

# Nueva Zelanda
Set-Row 71 @("Nueva Zelanda", 1440, 9, 974, 454, 14, 0, 12)

# Jamaica moves above Islas Feroe (rows 125 / 126 swap position)
Set-Row 125 @("Jamaica", 196, 23, 27, 164, 0, 0, 5)
Set-Row 126 @("Islas Feroe", 185, 0, 176, 9, 0, 0, 0)

# Siria moves above Mozambique (rows 166 / 167 swap position)
Set-Row 166 @("Siria", 39, 0, 5, 31, 0, 0, 3)
Set-Row 167 @("Mozambique", 39, 0, 8, 31, 0, 0, 0)

# Butan moves above Burundi (rows 209 / 210 swap position)
Set-Row 209 @("Butan", 5, 0, 2, 3, 0, 0, 0)
Set-Row 210 @("Burundi", 5, 0, 4, 0, 0, 0, 1)
